$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '88.894.71'
$ws.Cells.Item(2, 5).Value = '  +9.59%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.370.56'
$ws.Cells.Item(3, 5).Value = '  +7.42%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.04%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '223.30'
$ws.Cells.Item(5, 5).Value = '  +7.78%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '651.98'
$ws.Cells.Item(6, 5).Value = '  +5.51%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.414'
$ws.Cells.Item(7, 5).Value = '  +46.28%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.01%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.646'
$ws.Cells.Item(9, 5).Value = '  +11.41%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '3.368.51'
$ws.Cells.Item(10, 5).Value = '  +7.39%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.633'
$ws.Cells.Item(11, 5).Value = '  +8.95%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.0000288'
$ws.Cells.Item(12, 5).Value = '  +15.12%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '37.04'
$ws.Cells.Item(13, 5).Value = '  +18.52%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +2.62%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.992.23'
$ws.Cells.Item(15, 5).Value = '  +7.58%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '5.60'
$ws.Cells.Item(16, 5).Value = '  +6.31%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '88.766.02'
$ws.Cells.Item(17, 5).Value = '  +9.54%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.368.07'
$ws.Cells.Item(18, 5).Value = '  +7.62%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '14.99'
$ws.Cells.Item(19, 5).Value = '  +8.52%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -0.12%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '9.70'
$ws.Cells.Item(21, 5).Value = '  +9.08%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '460.08'
$ws.Cells.Item(22, 5).Value = '  +7.30%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '5.70'
$ws.Cells.Item(23, 5).Value = '  +12.19%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '7.47'
$ws.Cells.Item(24, 5).Value = '  +3.43%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '5.63'
$ws.Cells.Item(25, 5).Value = '  +8.02%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '12.86'
$ws.Cells.Item(26, 5).Value = '  +19.92%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '3.584.43'
$ws.Cells.Item(27, 5).Value = '  +8.51%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Litecoin'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '80.68'
$ws.Cells.Item(28, 5).Value = '  +6.74%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'PEPE'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.0000143'
$ws.Cells.Item(29, 5).Value = '  +18.86%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'Dai'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '1.00'
$ws.Cells.Item(30, 5).Value = '  -0.23%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Cronos'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.196'
$ws.Cells.Item(31, 5).Value = '  +42.84%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '9.44'
$ws.Cells.Item(32, 5).Value = '  +5.91%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '592.47'
$ws.Cells.Item(33, 5).Value = '  +1.80%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.00'
$ws.Cells.Item(34, 5).Value = '  +0.20%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +4.99%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.13'
$ws.Cells.Item(36, 5).Value = '  +7.87%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '7.26'
$ws.Cells.Item(37, 5).Value = '  +20.44%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -5.76%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '23.73'
$ws.Cells.Item(39, 5).Value = '  +5.11%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.437'
$ws.Cells.Item(40, 5).Value = '  +7.67%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +6.94%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '21.83'
$ws.Cells.Item(42, 5).Value = '  +5.43%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -0.07%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '3.19'
$ws.Cells.Item(44, 5).Value = '  +5.59%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '158.49'
$ws.Cells.Item(45, 5).Value = '  +0.03%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'USDe'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '1.00'
$ws.Cells.Item(46, 5).Value = '  +0.06%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'ImmutableX'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.45'
$ws.Cells.Item(47, 5).Value = '  +10.50%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '189.37'
$ws.Cells.Item(48, 5).Value = '  +2.12%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '46.59'
$ws.Cells.Item(49, 5).Value = '  +3.57%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '4.54'
$ws.Cells.Item(50, 5).Value = '  +8.90%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'ARBITRUM'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.671'
$ws.Cells.Item(51, 5).Value = '  +7.95%  '
